# Api Testing Test Cases Updated
#
# Adds a "1st recr" round of interview data (Arpit Dadhich) to the
# "karthik" sheet and moves the active sheet/selection from the
# "harsha" sheet over to "karthik".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "karthik"

# --- Header row gets two new (empty, but bold-styled) trailing columns ---
$ws.Range("E1:F1").Font.Bold = $true

# --- New row 3: first-round recruitment info ---
$ws.Range("A3").Value = "1st recr"
$ws.Range("B3").Value = "Frontend"
$ws.Range("C3").Value = """2022-03-31"""
$ws.Range("C3").NumberFormat = "h:mm"
$ws.Range("D3").Value = """2022-03-31"""
$ws.Range("D3").NumberFormat = "h:mm"

# --- New row 4: candidate details ---
$ws.Range("A4").Value = "Arpit Dadhich"
$ws.Range("B4").Value = "arpdadhich@deloitte.com"
$ws.Range("C4").Value = """9765432177"""
$ws.Range("D4").Value = "Not Eligible"
$ws.Range("F4").Value = "No"
$ws.Range("E4").Value = "Not Selected"

# --- Column widths: widen Name/Email column, add width for new column E ---
$ws.Columns.Item(2).ColumnWidth = 22.33
$ws.Columns.Item(5).ColumnWidth = 10

# --- Move the active sheet / selection from "harsha" to "karthik" ---
$ws.Activate() | Out-Null
$ws.Range("E4").Select() | Out-Null
